# The edit swaps the data of row 4 and row 5 for the columns that differ
# between the two otherwise-identical observation records (A, I, J, Q, R,
# X, Z, AB). Columns that already hold the same value in both rows are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Id (numeric) ---
$ws.Range("A4").Value = 131106008
$ws.Range("A5").Value = 131106010

# --- Antal / Enhet (text columns - use Copy so the "number-looking"
#     text in I stays text instead of being re-interpreted as numeric) ---
$ws.Range("I5").Copy($ws.Range("I4"))
$ws.Range("J5").Copy($ws.Range("J4"))
$ws.Range("I5").ClearContents()
$ws.Range("J5").ClearContents()

# --- Ost / Nord coordinates (numeric) ---
$ws.Range("Q4").Value = 612048
$ws.Range("R4").Value = 6945825
$ws.Range("Q5").Value = 612376
$ws.Range("R5").Value = 6945396

# --- Externid (text) ---
$ws.Range("X4").Value = "2025_1182"
$ws.Range("X5").Value = "2025_1180"

# --- Starttid / Sluttid (text) ---
$ws.Range("Z4").Value = "08:51"
$ws.Range("AB4").Value = "08:51"
$ws.Range("Z5").Value = "08:22"
$ws.Range("AB5").Value = "08:22"
